$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Pre-seed the new shared strings in the exact order they were first
# introduced in the authored workbook: Class C, Class D, col total, row total
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = "Class C"
$ws.Range("F20").Value = "Class D"
$ws.Range("F11").Value = "col total"
$ws.Range("B15").Value = "row total"

# ---------------------------------------------------------------------------
# Section 2: 3-class confusion matrix (rows 10-15)
# ---------------------------------------------------------------------------
$ws.Range("C10").Value = "TRUTH"

$ws.Range("C11").Value = "Class A"
$ws.Range("D11").Value = "Class B"

$ws.Range("A12").Value = "PREDICT"
$ws.Range("B12").Value = "Class A"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Formula = "=SUM(C12:E12)"
$ws.Range("H12").Value = "Obs. Acc:"
$ws.Range("I12").Formula = "=(C12+D13+E14)/F15"

$ws.Range("B13").Value = "Class B"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 18
$ws.Range("E13").Value = 10
$ws.Range("F13").Formula = "=SUM(C13:E13)"
$ws.Range("H13").Value = "Exp. Acc:"
$ws.Range("I13").Formula = "=(C15*F12 + D15*F13+E15*F14)/(F15^2)"

$ws.Range("B14").Value = "Class C"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 13
$ws.Range("E14").Value = 7
$ws.Range("F14").Formula = "=SUM(C14:E14)"
$ws.Range("H14").Value = "Kappa:"
$ws.Range("I14").Formula = "=(I12-I13)/(1-I13)"

$ws.Range("C15").Formula = "=SUM(C12:C14)"
$ws.Range("D15").Formula = "=SUM(D12:D14)"
$ws.Range("E15").Formula = "=SUM(E12:E14)"
$ws.Range("F15").Formula = '=IF(SUM(F12:F14)=SUM(C15:E15),SUM(C15:E15),CONCATENATE("R=",SUM(C15:E15),", C=",SUM(F12:F14)))'

# ---------------------------------------------------------------------------
# Section 3: 4-class confusion matrix (rows 19-25)
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = "TRUTH"

$ws.Range("C20").Value = "Class A"
$ws.Range("D20").Value = "Class B"
$ws.Range("E20").Value = "Class C"
$ws.Range("G20").Value = "col total"

$ws.Range("A21").Value = "PREDICT"
$ws.Range("B21").Value = "Class A"
$ws.Range("G21").Formula = "=SUM(C21:F21)"
$ws.Range("I21").Value = "Obs. Acc:"
$ws.Range("J21").Formula = "=(C21+D22+E23+F24)/G25"

$ws.Range("B22").Value = "Class B"
$ws.Range("D22").Value = 9
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22:G24").Formula = "=SUM(C22:F22)"
$ws.Range("I22").Value = "Exp. Acc:"
$ws.Range("J22").Formula = "=(C25*G21+D25*G22+E25*G23+F25*G24)/(G25^2)"

$ws.Range("B23").Value = "Class C"
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 2
$ws.Range("I23").Value = "Kappa:"
$ws.Range("J23").Formula = "=(J21-J22)/(1-J22)"

$ws.Range("B24").Value = "Class D"
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 6

$ws.Range("B25").Value = "row total"
$ws.Range("C25").Formula = "=SUM(C21:C24)"
$ws.Range("D25:F25").Formula = "=SUM(D21:D24)"
$ws.Range("G25").Formula = '=IF(SUM(G21:G24)=SUM(C25:F25),SUM(C25:F25),CONCATENATE("R=",SUM(C25:F25),", C=",SUM(G21:G24)))'

# ---------------------------------------------------------------------------
# Selection / view update to match final state
# ---------------------------------------------------------------------------
$ws.Range("I13").Select()
